$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1 - copy style from existing header cell H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data values for column I (I0) and J (IF)
$iValues = @(8, 8, 7, 5, 8, 5, 6, 6, 5, 8, 2)
$jValues = @(9, 9, 9, 5, 8, 5, 6, 6, 5, 8, 2)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
